$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Move the existing "26/06/2018" column (K) to its new home (Q),
#    freeing up K for the new "Wriiten" column.
# ---------------------------------------------------------------------
$ws.Range("Q1:Q15").Value = $ws.Range("K1:K15").Value
$ws.Range("K1:K15").ClearContents()

# ---------------------------------------------------------------------
# 2. New header row cells: J1=evidence, K1=Wriiten, L1=Total
# ---------------------------------------------------------------------
$ws.Range("J1").Value = "evidence"
$ws.Range("K1").Value = "Wriiten"
$ws.Range("L1").Value = "Total"

# New date columns R1 / S1 - copy the date number format used by C1:G1
$ws.Range("C1:D1").Copy()
$ws.Range("R1:S1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("R1").Value = 43227
$ws.Range("S1").Value = 43288

# ---------------------------------------------------------------------
# 3. Per-student data: evidence (J), written (K), total formula (L),
#    and the two new mark columns (R, S).
# ---------------------------------------------------------------------
$ws.Range("J2").Value = 42
$ws.Range("K2").Value = 74
$ws.Range("R2").Value = 58
$ws.Range("S2").Value = 16

$ws.Range("J3").Value = 44
$ws.Range("K3").Value = 70
$ws.Range("R3").Value = 50
$ws.Range("S3").Value = 15

$ws.Range("J4").Value = 44
$ws.Range("K4").Value = 70
$ws.Range("R4").Value = 62
$ws.Range("S4").Value = 17

$ws.Range("J5").Value = 48
$ws.Range("K5").Value = 78
$ws.Range("R5").Value = 64
$ws.Range("S5").Value = 18

$ws.Range("J6").Value = 48
$ws.Range("K6").Value = 79
$ws.Range("R6").Value = 60
$ws.Range("S6").Value = 18

$ws.Range("J7").Value = 44
$ws.Range("K7").Value = 70
$ws.Range("R7").Value = 60
$ws.Range("S7").Value = 17

$ws.Range("K8").Value = 58
$ws.Range("Q8").Value = 15

$ws.Range("J9").Value = 42
$ws.Range("K9").Value = 84
$ws.Range("R9").Value = 66
$ws.Range("S9").Value = 18

$ws.Range("J10").Value = 44
$ws.Range("K10").Value = 78
$ws.Range("R10").Value = 64
$ws.Range("S10").Value = 18

$ws.Range("J11").Value = 35
$ws.Range("R11").Value = 58
$ws.Range("S11").Value = 16

$ws.Range("J12").Value = 45
$ws.Range("K12").Value = 72
$ws.Range("R12").Value = 60
$ws.Range("S12").Value = 16

$ws.Range("J13").Value = 48
$ws.Range("K13").Value = 76
$ws.Range("R13").Value = 58
$ws.Range("S13").Value = 16

$ws.Range("J14").Value = 48
$ws.Range("K14").Value = 90
$ws.Range("R14").Value = 48
$ws.Range("S14").Value = 17

$ws.Range("J15").Value = 42
$ws.Range("R15").Value = 0

# ---------------------------------------------------------------------
# 4. Total column formulas (row 10/11 total D:K, everyone else C:K
#    because rows 10/11 have a non-numeric "C" cell).
# ---------------------------------------------------------------------
$ws.Range("L2").Formula = "=SUM(C2:K2)"
$ws.Range("L3").Formula = "=SUM(C3:K3)"
$ws.Range("L4").Formula = "=SUM(C4:K4)"
$ws.Range("L5").Formula = "=SUM(C5:K5)"
$ws.Range("L6").Formula = "=SUM(C6:K6)"
$ws.Range("L7").Formula = "=SUM(C7:K7)"
$ws.Range("L8").Formula = "=SUM(C8:K8)"
$ws.Range("L9").Formula = "=SUM(C9:K9)"
$ws.Range("L10").Formula = "=SUM(D10:K10)"
$ws.Range("L11").Formula = "=SUM(D11:K11)"
$ws.Range("L12").Formula = "=SUM(C12:K12)"
$ws.Range("L13").Formula = "=SUM(C13:K13)"
$ws.Range("L14").Formula = "=SUM(C14:K14)"
$ws.Range("L15").Formula = "=SUM(C15:K15)"

# ---------------------------------------------------------------------
# 5. Highlighting - whole rows 6,7,10 (B:L) plus L5 get a yellow fill;
#    the top result (L14) gets a yellow fill with red bold-ish font.
# ---------------------------------------------------------------------
$ws.Range("L14").Interior.Color = 65535
$ws.Range("L14").Font.Color = 255

$ws.Range("L5").Interior.Color = 65535
$ws.Range("B6:L7").Interior.Color = 65535
$ws.Range("B10:L10").Interior.Color = 65535

# ---------------------------------------------------------------------
# 6. Sheet view / selection / page setup cosmetics.
# ---------------------------------------------------------------------
$ws.Range("K11").Select()
$ws.PageSetup.Orientation = 1
